# Finish the 債務 (debt) sheet: fix the header row (B1:G1) to use proper
# field-name labels instead of a stray copy of row 2's data, and extend
# every row with the extra metadata columns H:N (property_category,
# category, date, legislator_name, legislator_id, source_file, index),
# matching the pattern already used on the other sheets (land/building/
# car/deposit/stock).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# ---- header row (row 1) ----
$ws.Cells.Item(1,2).Value = "species"
$ws.Cells.Item(1,3).Value = "debtor"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"
$ws.Cells.Item(1,6).Value = "register_date"
$ws.Cells.Item(1,7).Value = "register_reason"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = 11
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = 13

# ---- data rows (rows 2-8) ----
# columns: A index(year) B species C debtor D owner E total F register_date
#          G register_reason H property_category I category J date
#          K legislator_name L legislator_id M source_file N index

$rows = @(
    @{ A=95;  B="貸款";     C="廖國棟"; D="土地銀行台東分行臺東縣台東市中華路";   E=2529473; F="98年05月06日";  G="房貸";     N=95 },
    @{ A=96;  B="貸款";     C="陳姿伶"; D="華南銀行台東分行臺東縣台東市中華路";   E=4579121; F="90年01月12H";  G="房貸";     N=96 },
    @{ A=97;  B="貸款";     C="陳姿伶"; D="台灣企銀台東分行臺東縣台東市中華路";   E=6004284; F="94年10月11曰"; G="房貸";     N=97 },
    @{ A=98;  B="車貸";     C="陳姿伶"; D="台新銀行三重分行新北市三重區正義北路"; E=305434;  F="99年08月27日"; G="車貸";     N=98 },
    @{ A=99;  B="一般貸款"; C="廖國棟"; D="台東縣都蘭農會臺東縣東河鄉都蘭村都蘭"; E=4742000; F="93年08月31日"; G="—般貸款"; N=99 },
    @{ A=100; B="貸款";     C="廖國棟"; D="土地銀行台東分行臺東縣台東市中華路";   E=1848023; F="100年07月01日";G="信貸";     N=100 },
    @{ A=101; B="一般貸款"; C="陳姿伶"; D="台東縣都蘭農會臺東縣東河鄉都蘭村都蘭"; E=330750;  F="96年03月13曰"; G="—般貸款"; N=101 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value  = $row.A
    $ws.Cells.Item($r,2).Value  = $row.B
    $ws.Cells.Item($r,3).Value  = $row.C
    $ws.Cells.Item($r,4).Value  = $row.D
    $ws.Cells.Item($r,5).Value  = $row.E
    $ws.Cells.Item($r,6).Value  = $row.F
    $ws.Cells.Item($r,7).Value  = $row.G
    $ws.Cells.Item($r,8).Value  = "debt"
    $ws.Cells.Item($r,9).Value  = "normal"
    $ws.Cells.Item($r,10).Value = "2012-03-06"
    $ws.Cells.Item($r,11).Value = "廖國棟"
    $ws.Cells.Item($r,12).Value = 962
    $ws.Cells.Item($r,13).Value = "tmpec731"
    $ws.Cells.Item($r,14).Value = $row.N
    $r++
}
